$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'298.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-2.53%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'31.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-1.68%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.016"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-6.08%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07518"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.86%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.792"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.52%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.716"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'8.36%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.797"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'2.63%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9247"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.56%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1704"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.93%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07331"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-4.08%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07939"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.54%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03038"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-1.49%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09900"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.34%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001489"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-2.37%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.04648"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'2.04%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.006334"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.73%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.451"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.64%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'2.216"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-1.13%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.3291"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.72%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1336"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'1.27%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.559"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'8.20%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.1551"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-4.75%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001218"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.22%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004423"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-2.30%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001400"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'19.83%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001843"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'5.85%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01671"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'1.37%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04552"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.87%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007068"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-4.90%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1328"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-2.63%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002061"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-8.72%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01287"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-6.37%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006079"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-1.97%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.930"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'1.97%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-5.79%"
$ws.Range("E47").Style = "Normal"
